$wb = $excel.ActiveWorkbook

# ===== Overview sheet =====
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "cea2a183-35d0-4fd3-87ee-517af1424132.md"
$ws1.Range("B2").Value = "e2e\cea2a183-35d0-4fd3-87ee-517af1424132.md"
foreach ($h in $ws1.Hyperlinks) {
  if ($h.Range.Address() -eq '$B$2') {
    $h.TextToDisplay = "e2e\cea2a183-35d0-4fd3-87ee-517af1424132.md"
  }
}
$ws1.Range("G2").Value = "2016-08-21 01:04:42"

# ===== zh-cn sheet =====
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "cea2a183-35d0-4fd3-87ee-517af1424132.md"
$ws2.Range("G2").Value = "cea2a183-35d0-4fd3-87ee-517af1424132.6a2b4edacf22f9497fe531567fda7347eb0806a7.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-21 01:04:38"
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"
foreach ($h in @($ws2.Hyperlinks)) {
  if ($h.Range.Address() -eq '$A$2') {
    $h.TextToDisplay = "cea2a183-35d0-4fd3-87ee-517af1424132.md"
  }
  if ($h.Range.Address() -eq '$I$2') {
    $h.Delete()
  }
}
$ws2.Range("I2").Value = ""
$ws2.Range("I2").Style = "Normal"
$ws2.Columns.Item(9).ColumnWidth = 17.8
$ws2.Columns.Item(10).ColumnWidth = 20.8

# ===== de-de sheet =====
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "cea2a183-35d0-4fd3-87ee-517af1424132.md"
$ws3.Range("G2").Value = "cea2a183-35d0-4fd3-87ee-517af1424132.6a2b4edacf22f9497fe531567fda7347eb0806a7.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-21 01:04:42"
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"
foreach ($h in @($ws3.Hyperlinks)) {
  if ($h.Range.Address() -eq '$A$2') {
    $h.TextToDisplay = "cea2a183-35d0-4fd3-87ee-517af1424132.md"
  }
  if ($h.Range.Address() -eq '$I$2') {
    $h.Delete()
  }
}
$ws3.Range("I2").Value = ""
$ws3.Range("I2").Style = "Normal"
$ws3.Columns.Item(9).ColumnWidth = 17.8
$ws3.Columns.Item(10).ColumnWidth = 20.8

"done"
